$d = $word.ActiveDocument

# Create the new "Table Caption" paragraph style, based on the existing
# "Caption - Table" style (whose styleId is "Caption-Table"). The
# "Caption - Table" style was not being applied correctly, so a new
# style is introduced and used instead.
$newStyle = $d.Styles.Add("TableCaption", 1)
$newStyle.BaseStyle = "Caption-Table"
$newStyle.NameLocal = "Table Caption"
$newStyle.QuickStyle = $true

# Re-point the paragraph that previously used "Caption - Table" to the
# new "Table Caption" style.
foreach ($p in $d.Paragraphs) {
    if ($p.Style.NameLocal -eq "Caption - Table") {
        $p.Style = "TableCaption"
    }
}
